$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.221.73"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.861.89"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7038"
$ws.Range("E5").Value = "  -1.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.49"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07797"
$ws.Range("E8").Value = "  -2.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3110"
$ws.Range("E9").Value = "  -0.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.24"
$ws.Range("E10").Value = "  -4.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08001"
$ws.Range("E11").Value = "  -4.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.853.86"
$ws.Range("E12").Value = "  -2.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "93.52"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.166"
$ws.Range("E14").Value = "  -1.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6954"
$ws.Range("E15").Value = "  -3.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.364"
$ws.Range("E16").Value = "  +0.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008306"
$ws.Range("E17").Value = "  -1.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.199.02"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "252.16"
$ws.Range("E19").Value = "  +4.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.135.04"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.08"
$ws.Range("E21").Value = "  -1.33%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.488"
$ws.Range("E23").Value = "  -4.48%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1554"
$ws.Range("E25").Value = "  -1.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.984"
$ws.Range("E26").Value = "  -0.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.28"
$ws.Range("E27").Value = "  -2.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.78"
$ws.Range("E28").Value = "  +1.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.499"
$ws.Range("E29").Value = "  -0.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.278"
$ws.Range("E30").Value = "  -3.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.265"
$ws.Range("E31").Value = "  -1.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.207"
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05257"
$ws.Range("E33").Value = "  -2.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.885"
$ws.Range("E34").Value = "  -3.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7433"
$ws.Range("E35").Value = "  -0.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.156"
$ws.Range("E36").Value = "  -2.39%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01862"
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.254.50"
$ws.Range("E39").Value = "  -2.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.740"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.278"
$ws.Range("E41").Value = "  -4.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8984"
$ws.Range("E42").Value = "  -3.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "110.85"
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.41"
$ws.Range("E44").Value = "  -2.99%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.034.31"
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5199"
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.779"
$ws.Range("E49").Value = "  -1.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.400"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4292"
$ws.Range("E51").Value = "  -2.51%  "
